# Daily attendance processing - 2026-01-17 17:03:54
# Normalize the "Recorded By" (column G) text so the literal "System"
# token is ordered after the other recorder(s) instead of leading them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
    elseif ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
